# Update factsheets with text edits from COMM
#
# The "No. of 990 Filers w/ Gov Grants" column (column B on most sheets,
# column A on the "Overall" sheet) was stored as numeric values. Convert
# those cells to text (matching the rest of the sheet, which stores every
# other figure - dollar amounts, percentages - as text already) and add
# the missing "Total" row to the County sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    # Force Excel to store the value as text (not as a number) even though
    # it looks numeric, the same way a user typing `'362` into a cell would.
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---- Overall sheet: A2 (No. of 990 Filers w/ Gov Grants) ----
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "362"

# ---- County sheet: B2:B24 (No. of 990 Filers w/ Gov Grants) + new Total row ----
$wsCounty = $wb.Worksheets.Item("County")
Set-TextValue $wsCounty.Range("B2") "34"
Set-TextValue $wsCounty.Range("B3") "3"
Set-TextValue $wsCounty.Range("B4") "15"
Set-TextValue $wsCounty.Range("B5") "8"
Set-TextValue $wsCounty.Range("B6") "8"
Set-TextValue $wsCounty.Range("B7") "5"
Set-TextValue $wsCounty.Range("B8") "30"
Set-TextValue $wsCounty.Range("B9") "6"
Set-TextValue $wsCounty.Range("B10") "5"
Set-TextValue $wsCounty.Range("B11") "6"
Set-TextValue $wsCounty.Range("B12") "53"
Set-TextValue $wsCounty.Range("B13") "8"
Set-TextValue $wsCounty.Range("B14") "38"
Set-TextValue $wsCounty.Range("B15") "1"
Set-TextValue $wsCounty.Range("B16") "26"
Set-TextValue $wsCounty.Range("B17") "4"
Set-TextValue $wsCounty.Range("B18") "31"
Set-TextValue $wsCounty.Range("B19") "7"
Set-TextValue $wsCounty.Range("B20") "12"
Set-TextValue $wsCounty.Range("B21") "45"
Set-TextValue $wsCounty.Range("B22") "6"
Set-TextValue $wsCounty.Range("B23") "9"
Set-TextValue $wsCounty.Range("B24") "2"

# New row 25: Total (all columns stored as text, same as the rest of the sheet)
$wsCounty.Range("A25").Value = "Total"
Set-TextValue $wsCounty.Range("B25") "362"
Set-TextValue $wsCounty.Range("C25") '$261,596,078'
Set-TextValue $wsCounty.Range("D25") "6.11%"
Set-TextValue $wsCounty.Range("E25") "-25.29%"
Set-TextValue $wsCounty.Range("F25") "68.78%"

# ---- Congressional District sheet: B2:B3 ----
$wsCD = $wb.Worksheets.Item("Congressional District")
Set-TextValue $wsCD.Range("B2") "362"
Set-TextValue $wsCD.Range("B3") "362"

# ---- Size sheet: B2:B8 ----
$wsSize = $wb.Worksheets.Item("Size")
Set-TextValue $wsSize.Range("B2") "151"
Set-TextValue $wsSize.Range("B3") "92"
Set-TextValue $wsSize.Range("B4") "68"
Set-TextValue $wsSize.Range("B5") "21"
Set-TextValue $wsSize.Range("B6") "9"
Set-TextValue $wsSize.Range("B7") "21"
Set-TextValue $wsSize.Range("B8") "362"

# ---- Subsector sheet: B2:B14 ----
$wsSub = $wb.Worksheets.Item("Subsector")
Set-TextValue $wsSub.Range("B2") "34"
Set-TextValue $wsSub.Range("B3") "25"
Set-TextValue $wsSub.Range("B4") "26"
Set-TextValue $wsSub.Range("B5") "31"
Set-TextValue $wsSub.Range("B6") "4"
Set-TextValue $wsSub.Range("B7") "116"
Set-TextValue $wsSub.Range("B8") "1"
Set-TextValue $wsSub.Range("B9") "1"
Set-TextValue $wsSub.Range("B10") "25"
Set-TextValue $wsSub.Range("B11") "3"
Set-TextValue $wsSub.Range("B12") "95"
Set-TextValue $wsSub.Range("B13") "1"
Set-TextValue $wsSub.Range("B14") "362"
